# PlayerPerformance_3968.xlsx update
#  1. Clear the stray empty cell B2 on "ODI Batting"
#  2. Add a new "ODI Batting Extra" sheet (after "ODI Bowling") with
#     MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. ODI Batting!B2 should be a genuinely empty cell (it was an empty
#    inline string before) -> ClearContents removes it from the sheet.
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()

# ---------------------------------------------------------------------
# 2. New sheet "ODI Batting Extra", placed after "ODI Bowling"
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowling)
$extra.Name = "ODI Batting Extra"

# Match the page margins used by every other sheet in the workbook
# (0.75in / 1in / 0.5in) instead of this engine's blank-sheet defaults.
$extra.PageSetup.LeftMargin = 54
$extra.PageSetup.RightMargin = 54
$extra.PageSetup.TopMargin = 72
$extra.PageSetup.BottomMargin = 72
$extra.PageSetup.HeaderMargin = 36
$extra.PageSetup.FooterMargin = 36

# Header row -------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Reuse the same header style ("ODI Batting"!A1) that the other sheets use
# (bold font, thin border, centre/top alignment) instead of re-creating it.
$batting.Range("A1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows ----------------------------------------------------------
# columns: MATCH_CODE(text), BATTING_POSITION(number|blank), NUM_4(text),
#          NUM_6(text), PERCENT_RUNS_OF_TOTAL(text), MAN_OF_MATCH(text)
$rows = @(
    @("3436", $null, $null,  $null, $null,    "NO"),
    @("3438", 4,     "1",    "0",   "5.23%",  "NO"),
    @("3440", 3,     "2",    "0",   "19.91%", "NO"),
    @("3442", $null, $null,  $null, $null,    "NO"),
    @("3444", 4,     "0",    "0",   "0.54%",  "NO"),
    @("3851", $null, $null,  $null, $null,    "NO"),
    @("4206", 2,     "0",    "0",   "1.68%",  "NO"),
    @("4207", 1,     "1",    "0",   "2.02%",  "NO")
)

$textCols = @(1, 3, 4, 5, 6)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    foreach ($col in $textCols) {
        $value = $row[$col - 1]
        if ($value -ne $null) {
            $cell = $extra.Cells.Item($r, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $value
        }
    }

    $battingPos = $row[1]
    if ($battingPos -ne $null) {
        $extra.Cells.Item($r, 2).Value = $battingPos
    }
}

# Leave the workbook focused the way it started (first tab active),
# rather than leaving the newly added sheet selected.
$extra.Range("A1").Select()
$wb.Worksheets.Item(1).Activate()
